$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Report Generated On" timestamp
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:00 AM"

# Update Total Billed Amount
$ws.Range("C8").Value = 2732.71

# Clear the Scope ID # value (G10)
$ws.Range("G10").Value = ""

# Update the Wednesday (07/16/2025) line items pricing
$ws.Range("H16").Value = 478.55
$ws.Range("H17").Value = 478.55
$ws.Range("H18").Value = 478.55
$ws.Range("H19").Value = 1435.65

# Update the Friday (07/18/2025) line items pricing
$ws.Range("H24").Value = 648.53
$ws.Range("H25").Value = 648.53
$ws.Range("H26").Value = 1297.06
